# Updated via Streamlit Approval System
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("V12").Value = 7001
$ws.Range("AI12").Value = "HOLD"
$ws.Range("AJ12").Value = "HOLD"
$ws.Range("AK12").Value = "testing purpose"
$ws.Range("AL12").Value = "testing purpos"
$ws.Range("AM12").Value = "testing purpose"
$ws.Range("AN12").Value = "testing purpose"
$ws.Range("AO12").Value = "testing purpose"

# Row 13
$ws.Range("AK13").Value = "te"
$ws.Range("AL13").Value = "te"
$ws.Range("AM13").Value = "te"
$ws.Range("AN13").Value = "te"
$ws.Range("AO13").Value = "te"

# Rows 18-24: AN column only
$ws.Range("AN18").Value = "te"
$ws.Range("AN19").Value = "te"
$ws.Range("AN20").Value = "te"
$ws.Range("AN21").Value = "te"
$ws.Range("AN22").Value = "te"
$ws.Range("AN23").Value = "te"
$ws.Range("AN24").Value = "te"

# Rows 25-26: AK column only
$ws.Range("AK25").Value = "te"
$ws.Range("AK26").Value = "te"

# Row 27
$ws.Range("AK27").Value = "te"
$ws.Range("AL27").Value = "te"
$ws.Range("AM27").Value = "te"
$ws.Range("AN27").Value = "te"
$ws.Range("AO27").Value = "te"
